$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10. This shifts the existing rows 10-41
# down to rows 11-42 (matching the target diff where the old row 10
# ("1111"/Salame) becomes row 11, and so on through the old row 41
# ("7791070000078"/Papel higiénico) becoming row 42).
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the new article
# ("Sal fina sin tacc" / Rinsal / 7791004000013).
$ws.Range("A10").Value = 7791004000013
$ws.Range("B10").Value = "Sal"
$ws.Range("C10").Value = "fina"
$ws.Range("D10").Value = "sin tacc"
$ws.Range("E10").Value = "Rinsal"
$ws.Range("F10").Value = 500
$ws.Range("G10").Value = "gr."
$ws.Range("H10").Value = "Bolsa"
$ws.Range("I10").Value = "Sales"
$ws.Range("J10").Value = "Argentina"
$ws.Range("K10").Value = 20
$ws.Range("L10").Value = $false
$ws.Range("M10").Value = $true
$ws.Range("N10").Value = "C:\VentaSoft\Imágenes de artículos\7791004000013.png"
$ws.Range("O10").Value = $true

# The blank row created by Insert() doesn't pick up the bordered /
# right-aligned / wrap-text format ("Normal_Artículos" style) that the
# rest of column O uses, so restore it explicitly on the new O10 cell
# to match the surrounding rows (O9 above / O11 below).
$o10 = $ws.Range("O10")
$o10.HorizontalAlignment = -4152
$o10.WrapText = $true
$o10.Font.Name = "Calibri"
$o10.Font.Size = 11
$o10.Borders.Color = 12632256
$o10.Borders.LineStyle = 1
